# Appends the new match row (Odisha FC vs Bengaluru FC) to the bottom of the sheet,
# mirroring the formatting of the previous data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 27
$row = 28

$ws.Cells.Item($row, 1).Value = 27
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45230.64583333334
$ws.Cells.Item($row, 6).Value = "Odisha FC"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Bengaluru FC"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 2.13
$ws.Cells.Item($row, 11).Value = "27/10/2023 16:42"
$ws.Cells.Item($row, 12).Value = 2.23
$ws.Cells.Item($row, 13).Value = "31/10/2023 15:29"
$ws.Cells.Item($row, 14).Value = 3.44
$ws.Cells.Item($row, 15).Value = "27/10/2023 16:42"
$ws.Cells.Item($row, 16).Value = 3.5
$ws.Cells.Item($row, 17).Value = "31/10/2023 15:29"
$ws.Cells.Item($row, 18).Value = 3.4
$ws.Cells.Item($row, 19).Value = "27/10/2023 16:42"
$ws.Cells.Item($row, 20).Value = 3.23
$ws.Cells.Item($row, 21).Value = "31/10/2023 15:29"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/odisha-fc-bengaluru-fc/zsF6G8gj/"

# Copy the formatting (font/border/alignment/number format) from the row above
# so the new row visually matches the rest of the table (bold/boxed index cell,
# date-time formatted match-date cell) without introducing new style entries.
for ($col = 1; $col -le 22; $col++) {
    $ws.Cells.Item($srcRow, $col).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
